$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row 2
$ws.Range("B2").Value = 10.96649452729874
$ws.Range("C2").Value = 8.924616226414072
$ws.Range("D2").Value = 3.711527837427679
$ws.Range("F2").Value = 21.98721886483627
$ws.Range("G2").Value = 25.6515833154779
$ws.Range("H2").Value = 12.88429218846789
$ws.Range("I2").Value = 18.05561681493204
$ws.Range("M2").Value = 19.96770435713959

# row 3
$ws.Range("B3").Value = 10.39402281140068
$ws.Range("C3").Value = 8.465441934064286
$ws.Range("D3").Value = 3.699710036139603
$ws.Range("F3").Value = 21.83331500171699
$ws.Range("G3").Value = 25.27871501027947
$ws.Range("H3").Value = 12.90569977159586
$ws.Range("I3").Value = 18.12450452700004
$ws.Range("M3").Value = 19.36146973047553

# row 4
$ws.Range("B4").Value = 10.02660016241512
$ws.Range("C4").Value = 8.168555259246487
$ws.Range("D4").Value = 3.692438087971636
$ws.Range("F4").Value = 21.74728996592142
$ws.Range("G4").Value = 25.06044986821346
$ws.Range("H4").Value = 12.92237395700974
$ws.Range("I4").Value = 18.17269342649961
$ws.Range("M4").Value = 18.98773513805637

# row 5
$ws.Range("B5").Value = 9.873030160161111
$ws.Range("C5").Value = 8.043879366997974
$ws.Range("D5").Value = 3.689472156329233
$ws.Range("F5").Value = 21.71439526781253
$ws.Range("G5").Value = 24.97431758463394
$ws.Range("H5").Value = 12.93005252557474
$ws.Range("I5").Value = 18.19380139418988
$ws.Range("M5").Value = 18.83531273522645

# row 6
$ws.Range("B6").Value = 9.847303220217041
$ws.Range("C6").Value = 8.022956282107641
$ws.Range("D6").Value = 3.688979552745403
$ws.Range("F6").Value = 21.7090644612668
$ws.Range("G6").Value = 24.96018881303102
$ws.Range("H6").Value = 12.93138079856091
$ws.Range("I6").Value = 18.19739485332407
$ws.Range("M6").Value = 18.81000343073894

# row 7
$ws.Range("B7").Value = 10.02454438576283
$ws.Range("C7").Value = 8.166888686596895
$ws.Range("D7").Value = 3.692398096850688
$ws.Range("F7").Value = 21.74683754886556
$ws.Range("G7").Value = 25.05927670797291
$ws.Range("H7").Value = 12.9224739401791
$ws.Range("I7").Value = 18.17297215606003
$ws.Range("M7").Value = 18.98567966504216

# row 8
$ws.Range("B8").Value = 10.77249920930579
$ws.Range("C8").Value = 8.769432607705852
$ws.Range("D8").Value = 3.707457176566879
$ws.Range("F8").Value = 21.93241423484964
$ws.Range("G8").Value = 25.52088130492521
$ws.Range("H8").Value = 12.89093865788066
$ws.Range("I8").Value = 18.07814004173171
$ws.Range("M8").Value = 19.75913610480661

# row 9
$ws.Range("B9").Value = 12.10715869002118
$ws.Range("C9").Value = 9.830053371377517
$ws.Range("D9").Value = 3.736801751344916
$ws.Range("F9").Value = 22.36197388727946
$ws.Range("G9").Value = 26.50481567308241
$ws.Range("H9").Value = 12.85726552282412
$ws.Range("I9").Value = 17.93939527637786
$ws.Range("M9").Value = 21.25343546792334

# row 10
$ws.Range("B10").Value = 13.00112180194837
$ws.Range("C10").Value = 10.5331659318604
$ws.Range("D10").Value = 3.758167058855791
$ws.Range("F10").Value = 22.71516395272717
$ws.Range("G10").Value = 27.26728442900366
$ws.Range("H10").Value = 12.8498935682508
$ws.Range("I10").Value = 17.86687446469074
$ws.Range("M10").Value = 22.32467994549208

# row 11
$ws.Range("B11").Value = 13.41557174894834
$ws.Range("C11").Value = 10.83616646087028
$ws.Range("D11").Value = 3.767827559015419
$ws.Range("F11").Value = 22.88340505359035
$ws.Range("G11").Value = 27.62089168449149
$ws.Range("H11").Value = 12.85034487271261
$ws.Range("I11").Value = 17.84039721752676
$ws.Range("M11").Value = 22.80385037365323

# row 12
$ws.Range("B12").Value = 13.57301556516476
$ws.Range("C12").Value = 10.9484622656532
$ws.Range("D12").Value = 3.771475862665064
$ws.Range("F12").Value = 22.94814782476
$ws.Range("G12").Value = 27.75561001005929
$ws.Range("H12").Value = 12.8510649452945
$ws.Range("I12").Value = 17.83131732652009
$ws.Range("M12").Value = 22.98395392938794

# row 13
$ws.Range("B13").Value = 13.5392745747652
$ws.Range("C13").Value = 10.92438632879497
$ws.Range("D13").Value = 3.770690599740027
$ws.Range("F13").Value = 22.93415925501764
$ws.Range("G13").Value = 27.72656223777377
$ws.Range("H13").Value = 12.85088541237946
$ws.Range("I13").Value = 17.83323059578647
$ws.Range("M13").Value = 22.94522798011678

# row 14
$ws.Range("B14").Value = 13.42860017292977
$ws.Range("C14").Value = 10.84545418612793
$ws.Range("D14").Value = 3.76812791436396
$ws.Range("F14").Value = 22.88871109400105
$ws.Range("G14").Value = 27.63195960040473
$ws.Range("H14").Value = 12.85039309716322
$ws.Range("I14").Value = 17.83963118554294
$ws.Range("M14").Value = 22.81869556471186

# row 15
$ws.Range("B15").Value = 13.36031871361055
$ws.Range("C15").Value = 10.79678716404638
$ws.Range("D15").Value = 3.766556860215297
$ws.Range("F15").Value = 22.86100568686805
$ws.Range("G15").Value = 27.57411426330448
$ws.Range("H15").Value = 12.85016311142408
$ws.Range("I15").Value = 17.84367527661669
$ws.Range("M15").Value = 22.7410103288005

# row 16
$ws.Range("B16").Value = 12.97543087880976
$ws.Range("C16").Value = 10.51302316755338
$ws.Range("D16").Value = 3.757534414870181
$ws.Range("F16").Value = 22.70431670885583
$ws.Range("G16").Value = 27.24429789432833
$ws.Range("H16").Value = 12.84994082229248
$ws.Range("I16").Value = 17.86873672901693
$ws.Range("M16").Value = 22.29318563344902

# row 17
$ws.Range("B17").Value = 12.74807135223272
$ws.Range("C17").Value = 10.33461100457926
$ws.Range("D17").Value = 3.75198335049334
$ws.Range("F17").Value = 22.61009413543604
$ws.Range("G17").Value = 27.04358482547037
$ws.Range("H17").Value = 12.85078056014227
$ws.Range("I17").Value = 17.8857864175791
$ws.Range("M17").Value = 22.01624235253716

# row 18
$ws.Range("B18").Value = 12.61544839136121
$ws.Range("C18").Value = 10.23040802984237
$ws.Range("D18").Value = 3.748785089203594
$ws.Range("F18").Value = 22.55661613237341
$ws.Range("G18").Value = 26.92878524870492
$ws.Range("H18").Value = 12.85162159095321
$ws.Range("I18").Value = 17.89620539985373
$ws.Range("M18").Value = 21.85619477050568

# row 19
$ws.Range("B19").Value = 12.57022848590544
$ws.Range("C19").Value = 10.19485519594833
$ws.Range("D19").Value = 3.747701324490656
$ws.Range("F19").Value = 22.53863407680703
$ws.Range("G19").Value = 26.89003170714618
$ws.Range("H19").Value = 12.85196777299492
$ws.Range("I19").Value = 17.89983796754558
$ws.Range("M19").Value = 21.80188104990728

# row 20
$ws.Range("B20").Value = 12.77246627820009
$ws.Range("C20").Value = 10.35376749672939
$ws.Range("D20").Value = 3.752574844114207
$ws.Range("F20").Value = 22.62005054561077
$ws.Range("G20").Value = 27.0648854309376
$ws.Range("H20").Value = 12.85065409707679
$ws.Range("I20").Value = 17.88390798929411
$ws.Range("M20").Value = 22.04580313493105

# row 21
$ws.Range("B21").Value = 13.46121006625981
$ws.Range("C21").Value = 10.86870493892398
$ws.Range("D21").Value = 3.768880918088623
$ws.Range("F21").Value = 22.90203272604407
$ws.Range("G21").Value = 27.65972581034671
$ws.Range("H21").Value = 12.85052278360965
$ws.Range("I21").Value = 17.83772541341661
$ws.Range("M21").Value = 22.85589910443187

# row 22
$ws.Range("B22").Value = 13.91248216574182
$ws.Range("C22").Value = 11.19099504579071
$ws.Range("D22").Value = 3.779479391081242
$ws.Range("F22").Value = 23.09232108845736
$ws.Range("G22").Value = 28.05316737425541
$ws.Range("H22").Value = 12.85363859507491
$ws.Range("I22").Value = 17.81306344928129
$ws.Range("M22").Value = 23.37741755367314

# row 23
$ws.Range("B23").Value = 13.67363478519339
$ws.Range("C23").Value = 11.02029238525818
$ws.Range("D23").Value = 3.773828624269739
$ws.Range("F23").Value = 22.99023087428234
$ws.Range("G23").Value = 27.84280328472813
$ws.Range("H23").Value = 12.85168211960582
$ws.Range("I23").Value = 17.82571766791825
$ws.Range("M23").Value = 23.09985242218171

# row 24
$ws.Range("B24").Value = 12.76144327973704
$ws.Range("C24").Value = 10.34511192656507
$ws.Range("D24").Value = 3.752307451167123
$ws.Range("F24").Value = 22.61554709355996
$ws.Range("G24").Value = 27.0552535681151
$ws.Range("H24").Value = 12.85071015530621
$ws.Range("I24").Value = 17.88475530600725
$ws.Range("M24").Value = 22.03244127437831

# row 25
$ws.Range("B25").Value = 11.76094423723461
$ws.Range("C25").Value = 9.556358495932123
$ws.Range("D25").Value = 3.728891551042114
$ws.Range("F25").Value = 22.23897343911783
$ws.Range("G25").Value = 26.23106148195651
$ws.Range("H25").Value = 12.86333767205942
$ws.Range("I25").Value = 17.97180791702608
$ws.Range("M25").Value = 20.85299243637401
